$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.73
$ws.Range("S2").Value = 1.41
$ws.Range("T2").Value = 2.62
$ws.Range("AC2").Value = 9
$ws.Range("AK2").Value = 13
$ws.Range("AQ2").Value = 101
$ws.Range("AS2").Value = 301
$ws.Range("AX2").Value = 3.6
